$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All changed text-like cells are forced to Text number format ("@") before
# assigning the value, so that numeric-looking strings (e.g. "607.28") are not
# silently converted to floating point numbers by Excel, matching the source
# workbook where these are stored as plain inline strings.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.779.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.97%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.660.56'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.50%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.28'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.24'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.62%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.590'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.30%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +9.33%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.405'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.52%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.88'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.11%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.71%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.96'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +7.25%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000195'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +15.67%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.140.20'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.57%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.479.91'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.70%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.661.18'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.55%  '

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.29%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.89'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '361.94'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.05%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.44'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +6.32%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.45'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.99%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.70'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.76%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.60'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.55%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +17.14%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.99%  '

$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.166'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.60%  '

$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.19'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.14%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +7.30%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '541.12'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.43%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.09%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.84'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.67%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.61'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.31%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.38'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.22%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.85%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.70'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.72%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.02'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.15%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '162.73'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.11%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.10%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.02%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.40'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.92%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '166.34'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.80%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.17'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.09%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +8.16%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0614'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.29%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.13'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.51%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.661'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.83%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0265'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.13%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0990'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.01%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.83'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.06%  '
